# Update "想去人数" (want-to-go count) values in column F for the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, matching the
# regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row -> new F value, identical for both sheets.
$updates = @{
    5  = 15803
    9  = 15469
    11 = 9075
    12 = 390
    14 = 1016
    15 = 104
    20 = 61
    21 = 561
    29 = 96
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# Rows 35/39 on "展览" correspond to rows 37/41 on "全部类型" (the two
# sheets have their rows offset by 2 starting at this point).
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F35").Value = 261
$wsExhibit.Range("F39").Value = 5575

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F37").Value = 261
$wsAll.Range("F41").Value = 5575
